$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Exact set of StickerID values that need the "video:" prefix added,
# per the commit's fix to the answer-video lookup logic.
$targets = @(
    "Mientras_tanto.mp4",
    "Mucho_mucho_mas_tarde.mp4",
    "Muchos_Meses_Despues.mp4",
    "Nuevo_Narrador.mp4",
    "Uhhhh.mp4",
    "Una_deuda_Con_la_sociedad_mas_tarde.mp4",
    "Unos_momentos_despues.mp4",
    "Varias_bromas_malas_despues.mp4",
    "2000_años_más_tarde.mp4"
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $targets -contains $val) {
        $cell.Value = "video:" + $val
    }
}
